# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-15
$kValues = @{
    2  = 1
    3  = 0
    4  = 4
    5  = 5
    6  = 6
    7  = 2
    8  = 2
    9  = 2
    10 = 3
    11 = 4
    12 = 2
    13 = 1
    14 = 3
    15 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
